# Weekly refresh of the "Vega Monumental Concepción - Berenjena" sheet.
# Rows 32..64 each shift down by one row (row r's data moves to row r+1),
# a brand-new data row is written into row 32, and the old row 64's data
# (pushed out the bottom) is preserved as the new row 65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 32
$lastRow  = 64

# Walk from the bottom up so each source row is read before it gets
# overwritten by the row above it being shifted down.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $target = $r + 1

    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $i = $ws.Cells.Item($r, 9).Value2
    $j = $ws.Cells.Item($r, 10).Value2
    $k = $ws.Cells.Item($r, 11).Value2
    $l = $ws.Cells.Item($r, 12).Value2
    $m = $ws.Cells.Item($r, 13).Value2
    $n = $ws.Cells.Item($r, 14).Value2
    $o = $ws.Cells.Item($r, 15).Value2
    $p = $ws.Cells.Item($r, 16).Value2
    $q = $ws.Cells.Item($r, 17).Value2
    $s = $ws.Cells.Item($r, 18).Value2

    $ws.Cells.Item($target, 1).Value  = $a
    $ws.Cells.Item($target, 2).Value  = $b
    $ws.Cells.Item($target, 3).Value  = $c
    $ws.Cells.Item($target, 4).Value  = $d
    $ws.Cells.Item($target, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($target, 5).Value  = $e
    $ws.Cells.Item($target, 6).Value  = $f
    $ws.Cells.Item($target, 7).Value  = $g
    $ws.Cells.Item($target, 8).Value  = $h
    $ws.Cells.Item($target, 9).Value  = $i
    $ws.Cells.Item($target, 10).Value = $j
    $ws.Cells.Item($target, 11).Value = $k
    $ws.Cells.Item($target, 12).Value = $l
    $ws.Cells.Item($target, 13).Value = $m
    $ws.Cells.Item($target, 14).Value = $n
    $ws.Cells.Item($target, 15).Value = $o
    $ws.Cells.Item($target, 16).Value = $p
    $ws.Cells.Item($target, 17).Value = $q
    $ws.Cells.Item($target, 18).Value = $s
}

# New data point written into the now-vacated first row of the block.
$ws.Cells.Item($firstRow, 4).Value  = 44586
$ws.Cells.Item($firstRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($firstRow, 10).Value = 18500
$ws.Cells.Item($firstRow, 11).Value = 8000
$ws.Cells.Item($firstRow, 12).Value = 9000
$ws.Cells.Item($firstRow, 13).Value = 8514
$ws.Cells.Item($firstRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($firstRow, 16).Value = 142
